$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 4387.2,              4741, 3890, 0.06059975624084472),
    @(1, 4189.9,              4512, 3577, 0.06596214771270752),
    @(2, 4850.6,              5145, 4607, 0.0634625752766927),
    @(3, 5112.033333333334,   5522, 4570, 0.06348052024841308),
    @(4, 4328.2,              4562, 3976, 0.06034879684448242),
    @(5, 4986,                5314, 4421, 0.06562550067901611),
    @(6, 4425.866666666667,   4815, 3915, 0.06372884114583334),
    @(7, 4518.7,              4878, 4143, 0.06548502445220947),
    @(8, 4680.5,              5155, 4096, 0.06481369336446126),
    @(9, 4780.966666666666,   5243, 4389, 0.05949246088663737)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}
